$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cells in row 2 and 3 from numeric values to formatted text labels
# (order matches how the shared strings table was built in the target workbook)
$ws.Range("B2").Value = "12,000 words"
$ws.Range("B3").Value = "10,000 words"
$ws.Range("C2").Value = "4,000 words"

# Add new row 4 for Journal of Politics
$ws.Range("A4").Value = "JOURNAL OF POLITICS"
$ws.Range("D4").Value = "https://www.journals.uchicago.edu/journals/jop/instruct"
$ws.Range("B4").Value = "35 pages"
$ws.Range("C4").Value = "10 pages"
$ws.Range("E4").Value = 20190903

# Style D4 like D2/D3 (hyperlink style) and add the hyperlink
$ws.Range("D4").Style = $ws.Range("D3").Style
$ws.Hyperlinks.Add($ws.Range("D4"), "https://www.journals.uchicago.edu/journals/jop/instruct")
$ws.Range("D4").Value = "https://www.journals.uchicago.edu/journals/jop/instruct"
$ws.Range("D4").Style = $ws.Range("D3").Style

# Update selection to match the saved view state
$ws.Range("E5").Select()
